$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values are stored as text (matching original inlineStr formatting)
$ws.Range("D2:D5").NumberFormat = "@"
$ws.Range("D7:D26").NumberFormat = "@"
$ws.Range("D38:D43").NumberFormat = "@"
$ws.Range("D45:D46").NumberFormat = "@"
$ws.Range("D48:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"
$ws.Range("E38:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '306.83'
$ws.Range("E2").Value = '-4.01%'
$ws.Range("G2").Value = '16'
$ws.Range("D3").Value = '39.91'
$ws.Range("E3").Value = '-6.34%'
$ws.Range("G3").Value = '16'
$ws.Range("D4").Value = '5.042'
$ws.Range("E4").Value = '-3.12%'
$ws.Range("G4").Value = '16'
$ws.Range("D5").Value = '0.07613'
$ws.Range("E5").Value = '-6.94%'
$ws.Range("G5").Value = '16'
$ws.Range("E6").Value = '-2.44%'
$ws.Range("G6").Value = '16'
$ws.Range("D7").Value = '1.596'
$ws.Range("E7").Value = '-9.85%'
$ws.Range("G7").Value = '16'
$ws.Range("D8").Value = '0.9100'
$ws.Range("E8").Value = '-3.75%'
$ws.Range("G8").Value = '16'
$ws.Range("D9").Value = '0.1032'
$ws.Range("E9").Value = '-8.53%'
$ws.Range("G9").Value = '16'
$ws.Range("D10").Value = '0.1752'
$ws.Range("E10").Value = '-6.18%'
$ws.Range("G10").Value = '16'
$ws.Range("D11").Value = '0.09486'
$ws.Range("E11").Value = '0.66%'
$ws.Range("G11").Value = '16'
$ws.Range("D12").Value = '0.04465'
$ws.Range("E12").Value = '-4.66%'
$ws.Range("G12").Value = '16'
$ws.Range("D13").Value = '0.1051'
$ws.Range("E13").Value = '-0.58%'
$ws.Range("G13").Value = '16'
$ws.Range("D14").Value = '0.001236'
$ws.Range("E14").Value = '-4.95%'
$ws.Range("G14").Value = '16'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005831'
$ws.Range("E15").Value = '-0.31%'
$ws.Range("G15").Value = '16'
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '0.007491'
$ws.Range("E16").Value = '2,405.95%'
$ws.Range("G16").Value = '16'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.357'
$ws.Range("E17").Value = '0.17%'
$ws.Range("G17").Value = '16'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.419'
$ws.Range("E18").Value = '-4.46%'
$ws.Range("G18").Value = '16'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3312'
$ws.Range("E19").Value = '-1.45%'
$ws.Range("G19").Value = '16'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '6.924'
$ws.Range("E20").Value = '-6.97%'
$ws.Range("G20").Value = '16'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '0.1361'
$ws.Range("E21").Value = '-2.16%'
$ws.Range("G21").Value = '16'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '0.2817'
$ws.Range("E22").Value = '10.16%'
$ws.Range("G22").Value = '16'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '0.04148'
$ws.Range("E23").Value = '-0.69%'
$ws.Range("G23").Value = '16'
$ws.Range("D24").Value = '0.001210'
$ws.Range("E24").Value = '-2.28%'
$ws.Range("G24").Value = '16'
$ws.Range("D25").Value = '0.004074'
$ws.Range("E25").Value = '-4.59%'
$ws.Range("G25").Value = '16'
$ws.Range("D26").Value = '0.0001306'
$ws.Range("E26").Value = '6.64%'
$ws.Range("G26").Value = '16'
$ws.Range("G27").Value = '16'
$ws.Range("G28").Value = '16'
$ws.Range("G29").Value = '16'
$ws.Range("G30").Value = '16'
$ws.Range("G31").Value = '16'
$ws.Range("G32").Value = '16'
$ws.Range("G33").Value = '16'
$ws.Range("G34").Value = '16'
$ws.Range("G35").Value = '16'
$ws.Range("G36").Value = '16'
$ws.Range("G37").Value = '16'
$ws.Range("D38").Value = '0.02448'
$ws.Range("E38").Value = '-7.05%'
$ws.Range("G38").Value = '16'
$ws.Range("D39").Value = '0.05148'
$ws.Range("E39").Value = '-8.37%'
$ws.Range("G39").Value = '16'
$ws.Range("D40").Value = '0.007942'
$ws.Range("E40").Value = '-2.49%'
$ws.Range("G40").Value = '16'
$ws.Range("D41").Value = '0.1302'
$ws.Range("E41").Value = '-7.08%'
$ws.Range("G41").Value = '16'
$ws.Range("D42").Value = '0.007096'
$ws.Range("E42").Value = '8.13%'
$ws.Range("G42").Value = '16'
$ws.Range("D43").Value = '0.001958'
$ws.Range("E43").Value = '-6.55%'
$ws.Range("G43").Value = '16'
$ws.Range("E44").Value = '9.97%'
$ws.Range("G44").Value = '16'
$ws.Range("D45").Value = '0.3055'
$ws.Range("E45").Value = '-12.26%'
$ws.Range("G45").Value = '16'
$ws.Range("D46").Value = '0.00006443'
$ws.Range("E46").Value = '-5.00%'
$ws.Range("G46").Value = '16'
$ws.Range("E47").Value = '0.09%'
$ws.Range("G47").Value = '16'
$ws.Range("D48").Value = '0.003013'
$ws.Range("E48").Value = '-26.76%'
$ws.Range("G48").Value = '16'
$ws.Range("D49").Value = '0.004584'
$ws.Range("E49").Value = '36.66%'
$ws.Range("G49").Value = '16'
$ws.Range("E50").Value = '0.09%'
$ws.Range("G50").Value = '16'
$ws.Range("D51").Value = '0.0002009'
$ws.Range("E51").Value = '0.09%'
$ws.Range("G51").Value = '16'
